# The commit adds one new weekly price-report row for "Acelga" (Vega
# Modelo de Temuco) that was missing from the consolidated sheet. It is
# inserted right before the existing row 167, which pushes every
# subsequent row (old 167..308) down by one (new 168..309).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 167, shifting rows 167-308 down to 168-309.
$ws.Range("A167").EntireRow.Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A167").Value = 10
$ws.Range("B167").Value = "Vega Modelo de Temuco"
$ws.Range("C167").Value = "La Araucanía"
$ws.Range("D167").Value = 44669
$ws.Range("E167").Value = 9
$ws.Range("F167").Value = 100112009
$ws.Range("G167").Value = "Acelga"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 35
$ws.Range("K167").Value = 9000
$ws.Range("L167").Value = 9000
$ws.Range("M167").Value = 9000
$ws.Range("N167").Value = "$/docena de atados (12 kilos)"
$ws.Range("O167").Value = "Provincia de Cautín"
$ws.Range("P167").Value = 750
$ws.Range("Q167").Value = 12
$ws.Range("R167").Value = "Hortaliza"
